# Helper: convert EMU (English Metric Units, as used in raw OOXML) to the
# points used by the Shape.Left/Top/Width/Height COM properties.
# A small sub-EMU epsilon is added because the host's pt->EMU conversion
# truncates rather than rounds, which can otherwise drop the target value
# by 1 EMU for values that aren't exactly representable as a double.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + (0.49 / 12700.0)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "ZoneTexte 24" currently shows "(7)" -> becomes "(8)" ---
$zt24 = $s.Shapes.Item("ZoneTexte 24")
$zt24.TextFrame.TextRange.Text = "(8)"

# --- 2) "ZoneTexte 26" currently shows "(8)" -> becomes "(9)", and its
#        left offset moves from 1580786 EMU to 1494086 EMU ---
$zt26 = $s.Shapes.Item("ZoneTexte 26")
$zt26.Left = EmuToPt 1494086
$zt26.TextFrame.TextRange.Text = "(9)"

# --- 3) "ZoneTexte 28" currently shows "(9)" -> becomes "(10)", and its
#        width grows from 561372 EMU to 732893 EMU ---
$zt28 = $s.Shapes.Item("ZoneTexte 28")
$zt28.Width = EmuToPt 732893
$zt28.TextFrame.TextRange.Text = "(10)"

# --- 4) New ellipse marker "Ellipse 41", cloned from the existing
#        "Ellipse 25" ellipse shape/style (same noFill/red outline look) ---
$ellipseSrc = $s.Shapes.Item("Ellipse 25")
$ellipse41 = $ellipseSrc.Duplicate()
$ellipse41.Name = "Ellipse 41"
$ellipse41.Left = EmuToPt 1680726
$ellipse41.Top = EmuToPt 4831159
$ellipse41.Width = EmuToPt 576064
$ellipse41.Height = EmuToPt 864096

# --- 5) New text label "ZoneTexte 44" with text "(7)", cloned from the
#        existing "ZoneTexte 24" text box (same structure/formatting) ---
$textSrc = $s.Shapes.Item("ZoneTexte 24")
$zoneTexte44 = $textSrc.Duplicate()
$zoneTexte44.Name = "ZoneTexte 44"
$zoneTexte44.Left = EmuToPt 1206054
$zoneTexte44.Top = EmuToPt 5479231
$zoneTexte44.Width = EmuToPt 561372
$zoneTexte44.Height = EmuToPt 461665
$zoneTexte44.TextFrame.TextRange.Text = "(7)"
